$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = '@'
$cell.Value = '59.059.64'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.69%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.732.25'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.65%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.01%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '508.71'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.17%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '142.07'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.12%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.27%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.535'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.56%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.742.86'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.33%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.12'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.38%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.77%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.37%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.65%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.206.77'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.76%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '58.968.77'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.09%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '21.88'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.17%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0000137'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.77%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.731.44'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.69%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.76'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.51%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.02'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.20%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '346.00'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.02%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.27'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.46%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.05%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.61'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.44%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '63.25'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.13%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 2)
$cell.NumberFormat = '@'
$cell.Value = 'Polygon'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.427'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.59%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 2)
$cell.NumberFormat = '@'
$cell.Value = 'Kaspa'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.173'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.36%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.996'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.23%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0843'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.81%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.52'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.12%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.13%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.62'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.04%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '19.21'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.93%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '149.52'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.06%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.20%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.42'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.77%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.959'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.33%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.61%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '36.15'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.90%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.40'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.36%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.55'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.87%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.184.73'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -6.40%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0560'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.21%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.994'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.43%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.605'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.82%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '19.12'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -7.78%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.81'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.50%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.17%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0228'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.61%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0888'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.11%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '18.17'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.02%  '
$cell.Style = 'Normal'
